# Update LTPP / NGCS sections info sheet:
# Append 19 new LTPP section rows (rows 7-25) below the existing 5 data rows,
# reusing the same "Section Data Type" / "State" values and sharing the same
# layout as the existing rows (columns A-I populated, J..Z left blank).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 - IH 30 (Dallas)
$ws.Range("A7").Value = "LTPP Sections"
$ws.Range("B7").Value = "48-5035-IH30-LTPP"
$ws.Range("C7").Value = "IH 30"
$ws.Range("E7").Value = "Texas"
$ws.Range("F7").Value = "Dallas"
$ws.Range("G7").Value = "Dallas"
$ws.Range("I7").Value = "32.798345, -96.681312"

# Row 8 - IH 10 (Yoakum / Gonzales)
$ws.Range("A8").Value = "LTPP Sections"
$ws.Range("B8").Value = "48-5154-IH10-LTPP"
$ws.Range("C8").Value = "IH 10"
$ws.Range("D8").Value = "MP 655.5"
$ws.Range("E8").Value = "Texas"
$ws.Range("F8").Value = "Yoakum"
$ws.Range("G8").Value = "Gonzales"
$ws.Range("H8").Value = "West"
$ws.Range("I8").Value = "-"

# Row 9 - IH 27 (Amarillo / Randall)
$ws.Range("A9").Value = "LTPP Sections"
$ws.Range("B9").Value = "48-5336-IH27-LTPP"
$ws.Range("C9").Value = "IH 27"
$ws.Range("E9").Value = "Texas"
$ws.Range("F9").Value = "Amarillo"
$ws.Range("G9").Value = "Randall"
$ws.Range("I9").Value = "34.969167, -101.871829"

# Row 10 - IH 40 (Amarillo / Carson)
$ws.Range("A10").Value = "LTPP Sections"
$ws.Range("B10").Value = "48-5335-IH40-LTPP"
$ws.Range("C10").Value = "IH 40"
$ws.Range("E10").Value = "Texas"
$ws.Range("F10").Value = "Amarillo"
$ws.Range("G10").Value = "Carson"
$ws.Range("I10").Value = "35.2107952,-101.1275216"

# Row 11 - Loop 289 (Lubbock)
$ws.Range("A11").Value = "LTPP Sections"
$ws.Range("B11").Value = "48-1111-Loop289-LTPP"
$ws.Range("C11").Value = "Loop 289"
$ws.Range("E11").Value = "Texas"
$ws.Range("F11").Value = "Lubbock"
$ws.Range("G11").Value = "Lubbock"
$ws.Range("I11").Value = "33.531515, -101.804841"

# Row 12 - US 54 (El Paso)
$ws.Range("A12").Value = "LTPP Sections"
$ws.Range("B12").Value = "48-3779-US54-LTPP"
$ws.Range("C12").Value = "US 54"
$ws.Range("E12").Value = "Texas"
$ws.Range("F12").Value = "El Paso"
$ws.Range("G12").Value = "El Paso"
$ws.Range("I12").Value = "31.790822, -106.440705"

# Row 13 - BI 20E (Odessa / Ector)
$ws.Range("A13").Value = "LTPP Sections"
$ws.Range("B13").Value = "48-5278-BI20E-LTPP"
$ws.Range("C13").Value = "BI 20E"
$ws.Range("E13").Value = "Texas"
$ws.Range("F13").Value = "Odessa"
$ws.Range("G13").Value = "Ector"
$ws.Range("I13").Value = "31.925035, -102.213078"

# Row 14 - US 90 (San Antonio / Bexar)
$ws.Range("A14").Value = "LTPP Sections"
$ws.Range("B14").Value = "48-1096-US90-LTPP"
$ws.Range("C14").Value = "US 90"
$ws.Range("E14").Value = "Texas"
$ws.Range("F14").Value = "San Antonio"
$ws.Range("G14").Value = "Bexar"
$ws.Range("I14").Value = "29.35514067,-98.83470903"

# Row 15 - IH 10 (Yoakum / Gonzales), second GPS point
$ws.Range("A15").Value = "LTPP Sections"
$ws.Range("B15").Value = "48-5154-IH10-LTPP"
$ws.Range("C15").Value = "IH 10"
$ws.Range("E15").Value = "Texas"
$ws.Range("F15").Value = "Yoakum"
$ws.Range("G15").Value = "Gonzales"
$ws.Range("I15").Value = "29.69240860,-97.23887494"

# Row 16 - SH 71 (Yoakum / Colorado)
$ws.Range("A16").Value = "LTPP Sections"
$ws.Range("B16").Value = "48-5024-SH71-LTPP"
$ws.Range("C16").Value = "SH 71"
$ws.Range("E16").Value = "Texas"
$ws.Range("F16").Value = "Yoakum"
$ws.Range("G16").Value = "Colorado"
$ws.Range("I16").Value = "29.73128995,-96.58131749"

# Row 17 - SH 146 (Houston / Harris)
$ws.Range("A17").Value = "LTPP Sections"
$ws.Range("B17").Value = "48-3010-SH146-LTPP"
$ws.Range("C17").Value = "SH 146"
$ws.Range("E17").Value = "Texas"
$ws.Range("F17").Value = "Houston"
$ws.Range("G17").Value = "Harris"
$ws.Range("I17").Value = "29.79191853,-94.90662801"

# Row 18 - FM 2223 (Bryan / Brazos)
$ws.Range("A18").Value = "LTPP Sections"
$ws.Range("B18").Value = "48-0802-FM2223-LTPP"
$ws.Range("C18").Value = "FM 2223"
$ws.Range("E18").Value = "Texas"
$ws.Range("F18").Value = "Bryan"
$ws.Range("G18").Value = "Brazos"
$ws.Range("I18").Value = "30.78784069,-96.41133312"

# Row 19 - SH 195 (Waco / Bell)
$ws.Range("A19").Value = "LTPP Sections"
$ws.Range("B19").Value = "48-A808-SH195-LTPP"
$ws.Range("C19").Value = "SH 195"
$ws.Range("E19").Value = "Texas"
$ws.Range("F19").Value = "Waco"
$ws.Range("G19").Value = "Bell"
$ws.Range("I19").Value = "30.98926033,-97.76215001"

# Row 20 - SH 121 (Fort Worth / Tarrant)
$ws.Range("A20").Value = "LTPP Sections"
$ws.Range("B20").Value = "48-5284-SH121-LTPP"
$ws.Range("C20").Value = "SH 121"
$ws.Range("E20").Value = "Texas"
$ws.Range("F20").Value = "Fort Worth"
$ws.Range("G20").Value = "Tarrant"
$ws.Range("I20").Value = "32.90914988,-97.09772845"

# Row 21 - SH 121 (Fort Worth / Tarrant), second segment
$ws.Range("A21").Value = "LTPP Sections"
$ws.Range("B21").Value = "48-5283-SH121-LTPP"
$ws.Range("C21").Value = "SH 121"
$ws.Range("E21").Value = "Texas"
$ws.Range("F21").Value = "Fort Worth"
$ws.Range("G21").Value = "Tarrant"
$ws.Range("I21").Value = "32.86419812,-97.10191135"

# Row 22 - US 287 (Fort Worth / Tarrant)
$ws.Range("A22").Value = "LTPP Sections"
$ws.Range("B22").Value = "48-5317-US287-LTPP"
$ws.Range("C22").Value = "US 287"
$ws.Range("E22").Value = "Texas"
$ws.Range("F22").Value = "Fort Worth"
$ws.Range("G22").Value = "Tarrant"
$ws.Range("I22").Value = "32.59556183,-97.14614634"

# Row 23 - IH 820 (Fort Worth / Tarrant)
$ws.Range("A23").Value = "LTPP Sections"
$ws.Range("B23").Value = "48-5301-IH820-LTPP"
$ws.Range("C23").Value = "IH 820"
$ws.Range("E23").Value = "Texas"
$ws.Range("F23").Value = "Fort Worth"
$ws.Range("G23").Value = "Tarrant"
$ws.Range("I23").Value = "32.71449651,-97.47855896"

# Row 24 - US 380 (Fort Worth / Wise)
$ws.Range("A24").Value = "LTPP Sections"
$ws.Range("B24").Value = "48-5310-US380-LTPP"
$ws.Range("C24").Value = "US 380"
$ws.Range("E24").Value = "Texas"
$ws.Range("F24").Value = "Fort Worth"
$ws.Range("G24").Value = "Wise"
$ws.Range("I24").Value = "33.23338830,-97.60971429"

# Row 25 - US 287 (Wichita Falls / Montague)
$ws.Range("A25").Value = "LTPP Sections"
$ws.Range("B25").Value = "48-5328-US287-LTPP"
$ws.Range("C25").Value = "US 287"
$ws.Range("E25").Value = "Texas"
$ws.Range("F25").Value = "Wichita Falls"
$ws.Range("G25").Value = "Montague"
$ws.Range("I25").Value = "33.58719157,-97.91291624"

# Match the saved cursor/selection position from the authored workbook.
[void]$ws.Range("N21").Select()
